# Refresh the "cryptos" price/volume table (Price in col D, Volume(1h) in col E).
# Price strings that look numeric (e.g. "1.005", "218.81") are written with a
# leading apostrophe so Excel stores them as text (matching the original
# inlineStr cells) instead of coercing them into floating-point numbers; the
# style is then reset to "Normal" so no stray number-format style sticks to
# the cell. Prices using dotted thousands separators (e.g. "26.170.59") are
# already non-numeric to Excel and need no special handling.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.170.59"
$ws.Range("E2").Value = "  -6.71%  "
$ws.Range("D3").Value = "1.674.07"
$ws.Range("E3").Value = "  -4.20%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'218.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.31%  "
$ws.Range("D6").Value = "'0.5079"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -12.50%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'0.2637"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").Value = "'0.06329"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.87%  "
$ws.Range("D10").Value = "'21.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.64%  "
$ws.Range("D11").Value = "'0.07391"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").Value = "1.676.22"
$ws.Range("E12").Value = "  -4.12%  "
$ws.Range("D13").Value = "'4.552"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").Value = "'0.5765"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.25%  "
$ws.Range("D15").Value = "1.897.54"
$ws.Range("E15").Value = "  -4.46%  "
$ws.Range("D16").Value = "'0.000008530"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "'64.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -12.19%  "
$ws.Range("D18").Value = "26.240.59"
$ws.Range("E18").Value = "  -6.55%  "
$ws.Range("D19").Value = "'4.964"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.59%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").Value = "'187.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.46%  "
$ws.Range("D23").Value = "'6.184"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.19%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").Value = "'143.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.82%  "
$ws.Range("D26").Value = "'7.632"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.69%  "
$ws.Range("D27").Value = "'0.1168"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("D28").Value = "'15.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").Value = "'1.309"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.39%  "
$ws.Range("D30").Value = "'0.05759"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.76%  "
$ws.Range("D31").Value = "'1.327"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.26%  "
$ws.Range("D32").Value = "'3.506"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.93%  "
$ws.Range("D33").Value = "'3.490"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.86%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "'1.006"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.83%  "
$ws.Range("D36").Value = "'0.5981"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.57%  "
$ws.Range("D37").Value = "'2.368"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").Value = "'2.634"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.74%  "
$ws.Range("D39").Value = "'0.01601"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.63%  "
$ws.Range("D40").Value = "1.085.32"
$ws.Range("E40").Value = "  -3.56%  "
$ws.Range("E41").Value = "  -6.02%  "
$ws.Range("D42").Value = "'0.8596"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").Value = "'1.004"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'99.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "1.820.00"
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("E46").Value = "  +3.59%  "
$ws.Range("D47").Value = "'56.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.07%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Value = "'8.064"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("D50").Value = "'0.4308"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("E51").Value = "  -3.52%  "
